# Apply text edits to slide 7 of the presentation:
#  1. "Rectangle 6" ASCII-art diagram: "Return Path Sub-TLVs" -> "Return Path Sub-TLV"
#     (keeping the box alignment by shifting the trailing space)
#  2. "Content Placeholder 2" heading: append clarifying phrase to
#     "Return Path TLV (value TBA2):"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Edit 1: the "Rectangle 6" shape holding the STAMP TLV diagram ---
$diagramShape = $s.Shapes.Item("Rectangle 6")
$diagramTextRange = $diagramShape.TextFrame.TextRange
$diagramPara = $diagramTextRange.Paragraphs(6)
$diagramPara.Runs(1).Text = "    |                    Return Path Sub-TLV                        |"

# --- Edit 2: the "Content Placeholder 2" shape holding the description ---
$descShape = $s.Shapes.Item("Content Placeholder 2")
$descTextRange = $descShape.TextFrame.TextRange
$descPara = $descTextRange.Paragraphs(1)
$descPara.Runs(1).Text = "Return Path TLV (value TBA2) to carry one sub-TLV for return path:"
